$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 4 new columns before column D (old D,E,F -> become H,I,J),
#    inheriting the formatting/width of column C, same as Excel does
#    when you select columns D:G and "Insert".
$ws.Columns("D:G").Insert()

# 2) Insert a new blank spacer row before the header row (old row 6 -> row 7)
$ws.Rows("6").Insert()

# 3) The title row (row 5) used to carry a thin bottom border; that border
#    now belongs to the new spacer row (row 6) instead. Copy row 5's
#    current format (font/fill/border) onto row 6 first...
$ws.Range("A5:J5").Copy()
$ws.Range("A6:J6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ...then strip the bottom border back off row 5.
$ws.Range("A5:J5").Borders.Item(9).LineStyle = -4142

# 4) Fill in the text for the 4 newly inserted header columns (D7:G7),
#    matching the style already used by the neighbouring header cells.
$ws.Range("C7").Copy()
$ws.Range("D7:G7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D7").Value = "Địa chỉ"
$ws.Range("E7").Value = "Quy mô"
$ws.Range("F7").Value = "Loại hình sản xuất"
$ws.Range("G7").Value = "Tên người đại diện"
